# Insert a new weekly price record for "Ajo" (Macroferia Regional de Talca)
# as row 316, pushing the existing row 316 (and everything after it) down
# by one row. This mirrors the OOXML diff: dimension grows from R384 to
# R385 and a brand-new record (date 44889) appears at row 316.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 316..384 down to 317..385, leaving a blank row 316 behind.
$ws.Range("A316").EntireRow.Insert()

# Populate the newly-inserted row 316 with the new record's data.
$ws.Range("A316").Value = 5
$ws.Range("B316").Value = "Macroferia Regional de Talca"
$ws.Range("C316").Value = "Maule"
$ws.Range("D316").Value = 44889
$ws.Range("E316").Value = 7
$ws.Range("F316").Value = 100112003
$ws.Range("G316").Value = "Ajo"
$ws.Range("H316").Value = "Chino"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 300
$ws.Range("K316").Value = 18000
$ws.Range("L316").Value = 18000
$ws.Range("M316").Value = 18000
$ws.Range("N316").Value = '$/malla 10 kilos'
$ws.Range("O316").Value = "China"
$ws.Range("P316").Value = 1800
$ws.Range("Q316").Value = 10
$ws.Range("R316").Value = "Hortaliza"
